$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.831.48"
$ws.Range("E2").Value = "  -0.89%  "

# Row 3
$ws.Range("D3").Value = "2.097.91"
$ws.Range("E3").Value = "  +2.39%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'245.63"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6
$ws.Range("D6").Value = "'0.654"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'54.70"
$ws.Range("E8").Value = "  -2.82%  "

# Row 9
$ws.Range("D9").Value = "'59.04"
$ws.Range("E9").Value = "  -1.69%  "

# Row 10
$ws.Range("D10").Value = "'0.368"
$ws.Range("E10").Value = "  -3.13%  "

# Row 11
$ws.Range("D11").Value = "'0.0764"
$ws.Range("E11").Value = "  -2.04%  "

# Row 12
$ws.Range("E12").Value = "  +1.35%  "

# Row 13
$ws.Range("D13").Value = "'0.917"
$ws.Range("E13").Value = "  +4.72%  "

# Row 14
$ws.Range("D14").Value = "'15.12"
$ws.Range("E14").Value = "  -6.72%  "

# Row 15
$ws.Range("D15").Value = "2.397.53"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16
$ws.Range("D16").Value = "'5.53"
$ws.Range("E16").Value = "  -3.00%  "

# Row 17
$ws.Range("D17").Value = "2.083.39"
$ws.Range("E17").Value = "  +1.54%  "

# Row 18
$ws.Range("D18").Value = "36.797.35"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
$ws.Range("D19").Value = "'17.19"
$ws.Range("E19").Value = "  -6.36%  "

# Row 20
$ws.Range("D20").Value = "'72.84"
$ws.Range("E20").Value = "  -2.23%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0882"
$ws.Range("E21").Value = "  -1.31%  "

# Row 22
$ws.Range("D22").Value = "'5.47"
$ws.Range("E22").Value = "  +1.68%  "

# Row 23
$ws.Range("D23").Value = "'238.92"
$ws.Range("E23").Value = "  +0.91%  "

# Row 24
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  -3.18%  "

# Row 26
$ws.Range("D26").Value = "'9.81"
$ws.Range("E26").Value = "  +2.98%  "

# Row 27
$ws.Range("D27").Value = "'2.16"
$ws.Range("E27").Value = "  -0.33%  "

# Row 28
$ws.Range("D28").Value = "'167.08"
$ws.Range("E28").Value = "  -1.40%  "

# Row 29
$ws.Range("D29").Value = "'20.89"
$ws.Range("E29").Value = "  +4.35%  "

# Row 30
$ws.Range("E30").Value = "  -1.20%  "

# Row 31
$ws.Range("D31").Value = "'5.25"
$ws.Range("E31").Value = "  +8.54%  "

# Row 32
$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  +3.81%  "

# Row 33
$ws.Range("D33").Value = "'4.68"
$ws.Range("E33").Value = "  +4.56%  "

# Row 34
$ws.Range("D34").Value = "'0.0610"
$ws.Range("E34").Value = "  -1.40%  "

# Row 35
$ws.Range("D35").Value = "'2.45"
$ws.Range("E35").Value = "  +9.52%  "

# Row 36
$ws.Range("E36").Value = "  -0.11%  "

# Row 37
$ws.Range("E37").Value = "  +3.83%  "

# Row 38
$ws.Range("E38").Value = "  -7.10%  "

# Row 39
$ws.Range("E39").Value = "  -4.78%  "

# Row 40
$ws.Range("E40").Value = "  +1.33%  "

# Row 41
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'4.91"
$ws.Range("E41").Value = "  -6.90%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0221"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43
$ws.Range("D43").Value = "'0.0956"
$ws.Range("E43").Value = "  -4.25%  "

# Row 44
$ws.Range("D44").Value = "'96.67"
$ws.Range("E44").Value = "  +1.12%  "

# Row 45
$ws.Range("D45").Value = "'2.84"
$ws.Range("E45").Value = "  -8.78%  "

# Row 46
$ws.Range("D46").Value = "1.408.66"
$ws.Range("E46").Value = "  +11.35%  "

# Row 47
$ws.Range("D47").Value = "'16.19"
$ws.Range("E47").Value = "  -6.11%  "

# Row 48
$ws.Range("D48").Value = "'7.57"
$ws.Range("E48").Value = "  +11.54%  "

# Row 49
$ws.Range("D49").Value = "'2.48"
$ws.Range("E49").Value = "  +2.10%  "

# Row 50
$ws.Range("D50").Value = "'2.91"
$ws.Range("E50").Value = "  +2.26%  "

# Row 51
$ws.Range("D51").Value = "2.286.37"
$ws.Range("E51").Value = "  +2.09%  "

